$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Chapter "2-Marco teorico" (row 3) has now been read -> "Leido" instead of "En proceso"
$ws.Range("C3").Value = "Leido"

# Chapter "4-Metodologia de Desarrollo" (row 5) is now in progress -> "En proceso"
$ws.Range("C5").Value = "En proceso"

# Move the active cell selection to C5
$ws.Range("C5").Select()
